# Update OverVoltage (OV) values for 4th and 5th week rows, and
# adjust the row heights for rows 1-37 from their previous values
# (16.5 / 18.75) to a uniform 18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row height adjustments (rows 1 through 37 -> 18) ---
$htRows = 1..37
foreach ($r in $htRows) {
    $ws.Rows.Item($r).RowHeight = 18
}

# --- OverVoltage (column F) updates: "4 4 4 4 3" -> "4 4 4 4 2.5" ---
$ovRows = @(4, 5, 9, 14, 15, 19, 22, 26, 27, 31, 36, 37, 42, 43)
foreach ($r in $ovRows) {
    $ws.Cells.Item($r, 6).Value = "4 4 4 4 2.5"
}
